$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row/column of the existing data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New header cells AD1:AF1 - copy formatting from the existing header style (AC1)
# so the new headers match the bold/centered/bordered header look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 72   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
